$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999204813972,
    0.99553801280933907,
    0.97742543838285245,
    0.96991336460522382,
    0.96288363594276938,
    0.94700577833399313,
    0.94414894941321414,
    0.93959180618295735,
    0.93621085820507877,
    0.93375481209997702,
    0.93339857891891986,
    0.93297689248913929,
    0.93387043447338036,
    0.93544760497305379,
    0.93800745504114647,
    0.94148170543665011,
    0.95157019721048797,
    0.95652768513969233,
    0.9939522330652566,
    0.98184276925998537,
    0.98044429434094149,
    0.97917979125045207,
    0.97039017938065053,
    0.9573688267384477,
    0.95091175849441623,
    0.94563594339254631,
    0.94397760797544172,
    0.9377106326693383,
    0.93376319251231743,
    0.93270629338903865,
    0.93961332768670958,
    0.94304843805648475,
    0.94363865973338423
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
